$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 16:52"

# Tenerife (row 34)
$ws.Range("B34").Value = 1391
$ws.Range("C34").Value = 512
$ws.Range("D34").Value = 798

# Gran Canaria (row 50)
$ws.Range("B50").Value = 500
$ws.Range("D50").Value = 232

# La Palma (row 56)
$ws.Range("B56").Value = 74
$ws.Range("C56").Value = 25
$ws.Range("D56").Value = 46

# Lanzarote (row 57)
$ws.Range("C57").Value = 17
$ws.Range("D57").Value = 49
